# Update KNN-imputed values in Sheet1 (terrestrial_mammals / BCE / 15 / seed1)
# to the refreshed algorithm output ("Update Name of Algo").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 7.645999999999999
$ws.Range("B6").Value = 6.637
$ws.Range("B7").Value = 5.281
$ws.Range("C7").Value = -13.383
$ws.Range("C12").Value = -11.087
$ws.Range("E13").Value = 16.684
$ws.Range("E14").Value = 17.135
$ws.Range("C15").Value = -13.385
$ws.Range("B16").Value = 5.532999999999999
$ws.Range("E16").Value = 16.544
$ws.Range("E19").Value = 16.54
$ws.Range("B20").Value = 7.794
$ws.Range("C20").Value = -11.992
$ws.Range("C21").Value = -12.679
$ws.Range("C22").Value = -12.692
$ws.Range("E22").Value = 16.66
$ws.Range("C23").Value = -12.451
$ws.Range("B28").Value = 6.296000000000001
$ws.Range("B29").Value = 5.209999999999999
$ws.Range("C29").Value = -11.358
$ws.Range("B32").Value = 6.256
$ws.Range("C34").Value = -12.44
$ws.Range("E36").Value = 16.499
$ws.Range("B40").Value = 9.263999999999999
$ws.Range("C42").Value = -11.999
$ws.Range("C43").Value = -13.867
$ws.Range("C44").Value = -13.94
$ws.Range("C45").Value = -13.376
$ws.Range("B46").Value = 4.944000000000001
$ws.Range("C46").Value = -14.281
$ws.Range("E46").Value = 16.618
$ws.Range("C50").Value = -13.933
$ws.Range("E50").Value = 16.487
$ws.Range("B51").Value = 5.059
$ws.Range("C51").Value = -12.057
$ws.Range("B52").Value = 5.459000000000001
$ws.Range("B57").Value = 5.763000000000001
$ws.Range("B59").Value = 5.243
$ws.Range("B62").Value = 5.902
$ws.Range("B66").Value = 4.961
$ws.Range("C66").Value = -10.897
$ws.Range("C67").Value = -11.354
$ws.Range("B73").Value = 7.316
$ws.Range("B74").Value = 8.944000000000001
$ws.Range("C79").Value = -12.062
$ws.Range("C84").Value = -13.68
$ws.Range("B92").Value = 6.425999999999999
$ws.Range("C92").Value = -10.955
$ws.Range("E95").Value = 17.518
$ws.Range("C97").Value = -11.87
$ws.Range("E97").Value = 17.114
$ws.Range("B100").Value = 6.135
